$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MoM")

# Insert a new column before column G ("transperadd"), shifting
# sess_mom..ECR_momper from G:P to H:Q.
$ws.Columns("G:G").Insert()

$ws.Cells.Item(1, 7).Value = "transperadd"
$ws.Cells.Item(2, 7).Value = 0.207643358689292
$ws.Cells.Item(3, 7).Value = 0.319885153283319

# Add two new trailing columns: transperadd_mom, transperadd_momper.
$ws.Cells.Item(1, 18).Value = "transperadd_mom"
$ws.Cells.Item(1, 19).Value = "transperadd_momper"
# Row 2 has no MoM value yet (nothing to compare against), so the cells
# stay present-but-empty like the rest of row 2's trailing columns.
$ws.Cells.Item(2, 18).Style = "Normal"
$ws.Cells.Item(2, 19).Style = "Normal"
$ws.Cells.Item(3, 18).Value = 0.112241794594027
$ws.Cells.Item(3, 19).Value = 0.540550852685739

# Column A now holds the "5/2013" / "6/2013" text labels instead of
# serial dates, so drop the date number format/style and set text values.
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(2, 1).Value = "5/2013"
$ws.Cells.Item(3, 1).Value = "6/2013"
